# Reorder the "Client Details" sheet's columns so that Project Name (old
# column C) moves to the front, pushing Account (old A) and Project ID
# (old B) one column to the right: A=Project Name, B=Account, C=Project ID.
# This is done the way a user would in Excel: cut column C and insert the
# cut cells before column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Client Details")

# Cut column C (Project Name) and insert it before column A.
$ws.Columns("C:C").Cut()
$ws.Columns("A:A").Insert()

# The AutoFilter range does not auto-track the column move, so turn it off
# and re-apply it over the new data range (now B1:C12 - Account/Project ID).
$ws.AutoFilterMode = $false
$rng = $ws.Range("B1:C12")
[void]$rng.AutoFilter()

# Re-sort by Project ID (now column C) ascending, same as before the move.
$rng.Sort($ws.Range("C1:C12"), 1, "", "", 1, "", 1, 1)

# Keep the workbook-level hidden _FilterDatabase name for this sheet in
# sync with the new filter range.
$wb.Names.Item("Client Details!_FilterDatabase").RefersTo = "='Client Details'!`$B`$1:`$C`$1"

# Final cursor position left by the edit.
[void]$ws.Range("A12").Select()
